$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates -----------------------------------------------------
# B5: "45 hours 37 minutes 23 seconds" -> "45 hours 29 minutes 28 seconds"
$ws.Range("B5").Value = "45 hours 29 minutes 28 seconds"

# G5: 894416550 -> 903373439
$ws.Range("G5").Value = 903373439

# U5: new cell "2937 seconds" (HSE (mt) column), right-aligned like the
# rest of the row's label cells (style index 2 in the original workbook)
$ws.Range("U5").Value = "2937 seconds"
$ws.Range("U5").HorizontalAlignment = -4152

# --- View state ---------------------------------------------------------
# Move the selection to U5 (also scroll so column M is the left-most
# visible column, matching the author's view when they added this cell)
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 13
$ws.Range("U5").Select()
